$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07479333333333334
$ws.Range("H2").Value = 0.22438
$ws.Range("I2").Value = 0.1617287198578621
$ws.Range("J2").Value = 0.1617287198578621
$ws.Range("M2").Value = 1.949849666666667
$ws.Range("N2").Value = 5.849549000000001
$ws.Range("O2").Value = 0.06676506732104066
$ws.Range("P2").Value = 0.06676506732104066
$ws.Range("Q2").Value = 0.1458357560688889
$ws.Range("R2").Value = 1.31252180462
$ws.Range("S2").Value = 0.01079782886905589
$ws.Range("T2").Value = 0.01079782886905589
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07479333333333334
$ws.Range("H3").Value = 0.22438
$ws.Range("I3").Value = 0.1617287198578621
$ws.Range("J3").Value = 0.1617287198578621
$ws.Range("O3").Value = 0.7967262871802238
$ws.Range("P3").Value = 0.7967262871802239
$ws.Range("Q3").Value = 1.740299008644445
$ws.Range("R3").Value = 15.6626910778
$ws.Range("S3").Value = 0.128853522502765
$ws.Range("T3").Value = 0.128853522502765
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07479333333333334
$ws.Range("H4").Value = 0.22438
$ws.Range("I4").Value = 0.1617287198578621
$ws.Range("J4").Value = 0.1617287198578621
$ws.Range("O4").Value = 0.1365086454987356
$ws.Range("P4").Value = 0.1365086454987356
$ws.Range("Q4").Value = 0.2981775099622222
$ws.Range("R4").Value = 2.68359758966
$ws.Range("S4").Value = 0.02207736848604122
$ws.Range("T4").Value = 0.02207736848604122
$ws.Range("I5").Value = 0.8382712801421379
$ws.Range("J5").Value = 0.8382712801421379
$ws.Range("M5").Value = 1.949849666666667
$ws.Range("N5").Value = 5.849549000000001
$ws.Range("O5").Value = 0.06676506732104066
$ws.Range("P5").Value = 0.06676506732104066
$ws.Range("Q5").Value = 0.7558949705272223
$ws.Range("R5").Value = 6.803054734745001
$ws.Range("S5").Value = 0.05596723845198478
$ws.Range("T5").Value = 0.05596723845198478
$ws.Range("I6").Value = 0.8382712801421379
$ws.Range("J6").Value = 0.8382712801421379
$ws.Range("O6").Value = 0.7967262871802238
$ws.Range("P6").Value = 0.7967262871802239
$ws.Range("S6").Value = 0.6678727646774588
$ws.Range("T6").Value = 0.6678727646774588
$ws.Range("I7").Value = 0.8382712801421379
$ws.Range("J7").Value = 0.8382712801421379
$ws.Range("O7").Value = 0.1365086454987356
$ws.Range("P7").Value = 0.1365086454987356
$ws.Range("S7").Value = 0.1144312770126944
$ws.Range("T7").Value = 0.1144312770126944